# Insert a new weekly price record for "Ajo" (garlic) at Terminal
# Hortofrutícola Agro Chillán. The new observation is inserted as row 62,
# pushing the existing rows 62-117 down to 63-118 (the sheet keeps its
# chronological/ordering convention of newest-first by inserting above the
# previous top record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 62..117 down to 63..118, leaving a blank row 62 to fill in.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A62").Value = 7
$ws.Range("B62").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C62").Value = "Ñuble"
$ws.Range("D62").Value = 44447
$ws.Range("E62").Value = 16
$ws.Range("F62").Value = 100112003
$ws.Range("G62").Value = "Ajo"
$ws.Range("H62").Value = "Chino"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 60
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 16000
$ws.Range("M62").Value = 15500
$ws.Range("N62").Value = "`$/caja 10 kilos"
$ws.Range("O62").Value = "China"
$ws.Range("P62").Value = 1550
$ws.Range("Q62").Value = 10
$ws.Range("R62").Value = "Hortaliza"
